# Generate Report for Handoff
# Adds two newly-discovered e2e files to the localization status report:
#   31a51f48-553e-4e6f-bf31-eec0f15bb255.md
#   ef7eeca8-c9e4-46c0-9841-5970a366c45a.md
# as new rows 4 and 5 on the Overview, zh-cn and de-de sheets, then grows
# each sheet's table to cover the new rows.
#
# Note: "True"/"False" and empty-string values must be entered with a
# leading apostrophe so Excel keeps them as literal text (matching the
# existing report data) instead of coercing to boolean / blank cells.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value2 = "31a51f48-553e-4e6f-bf31-eec0f15bb255.md"
$wsOverview.Range("B4").Value2 = "e2e\31a51f48-553e-4e6f-bf31-eec0f15bb255.md"
$wsOverview.Range("C4").Value2 = ".md"
$wsOverview.Range("D4").Value2 = "'"
$wsOverview.Range("E4").Value2 = "Ready for handoff"
$wsOverview.Range("F4").Value2 = "Ready for handoff"
$wsOverview.Range("G4").Value2 = "2016-08-12 04:43:20"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Range("A5").Value2 = "ef7eeca8-c9e4-46c0-9841-5970a366c45a.md"
$wsOverview.Range("B5").Value2 = "e2e\ef7eeca8-c9e4-46c0-9841-5970a366c45a.md"
$wsOverview.Range("C5").Value2 = ".md"
$wsOverview.Range("D5").Value2 = "'"
$wsOverview.Range("E5").Value2 = "Ready for handoff"
$wsOverview.Range("F5").Value2 = "Ready for handoff"
$wsOverview.Range("G5").Value2 = "2016-08-12 04:43:20"
$wsOverview.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/581f422deee33de35aa3257f53aa5bcd25473dcc/e2e/31a51f48-553e-4e6f-bf31-eec0f15bb255.md", "", "", "e2e\31a51f48-553e-4e6f-bf31-eec0f15bb255.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/581f422deee33de35aa3257f53aa5bcd25473dcc/e2e/ef7eeca8-c9e4-46c0-9841-5970a366c45a.md", "", "", "e2e\ef7eeca8-c9e4-46c0-9841-5970a366c45a.md")

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G5"))

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A4").Value2 = "31a51f48-553e-4e6f-bf31-eec0f15bb255.md"
$wsZhCn.Range("B4").Value2 = ".md"
$wsZhCn.Range("C4").Value2 = "Ready for handoff"
$wsZhCn.Range("D4").Value2 = "e2e"
$wsZhCn.Range("E4").Value2 = "ht"
$wsZhCn.Range("F4").Value2 = "'False"
$wsZhCn.Range("G4").Value2 = "31a51f48-553e-4e6f-bf31-eec0f15bb255.4314296d17028732499ce604519a357cf1aa06d4.zh-cn.xlf"
$wsZhCn.Range("H4").Value2 = "2016-08-12 04:43:15"
$wsZhCn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I4").Value2 = "'"
$wsZhCn.Range("J4").Value2 = "'"
$wsZhCn.Range("K4").Value2 = "0001-01-01 00:00:00"
$wsZhCn.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L4").Value2 = "'"
$wsZhCn.Range("M4").Value2 = "'True"
$wsZhCn.Range("N4").Value2 = "'"
$wsZhCn.Range("O4").Value2 = "'False"
$wsZhCn.Range("P4").Value2 = "'"

$wsZhCn.Range("A5").Value2 = "ef7eeca8-c9e4-46c0-9841-5970a366c45a.md"
$wsZhCn.Range("B5").Value2 = ".md"
$wsZhCn.Range("C5").Value2 = "Ready for handoff"
$wsZhCn.Range("D5").Value2 = "e2e"
$wsZhCn.Range("E5").Value2 = "ht"
$wsZhCn.Range("F5").Value2 = "'False"
$wsZhCn.Range("G5").Value2 = "ef7eeca8-c9e4-46c0-9841-5970a366c45a.8b6e6573c64d4bec93863aad80d8c624a38d0cef.zh-cn.xlf"
$wsZhCn.Range("H5").Value2 = "2016-08-12 04:43:15"
$wsZhCn.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I5").Value2 = "'"
$wsZhCn.Range("J5").Value2 = "'"
$wsZhCn.Range("K5").Value2 = "0001-01-01 00:00:00"
$wsZhCn.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L5").Value2 = "'"
$wsZhCn.Range("M5").Value2 = "'True"
$wsZhCn.Range("N5").Value2 = "'"
$wsZhCn.Range("O5").Value2 = "'False"
$wsZhCn.Range("P5").Value2 = "'"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/581f422deee33de35aa3257f53aa5bcd25473dcc/e2e/31a51f48-553e-4e6f-bf31-eec0f15bb255.md", "", "", "31a51f48-553e-4e6f-bf31-eec0f15bb255.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/581f422deee33de35aa3257f53aa5bcd25473dcc/e2e/ef7eeca8-c9e4-46c0-9841-5970a366c45a.md", "", "", "ef7eeca8-c9e4-46c0-9841-5970a366c45a.md")

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P5"))

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A4").Value2 = "31a51f48-553e-4e6f-bf31-eec0f15bb255.md"
$wsDeDe.Range("B4").Value2 = ".md"
$wsDeDe.Range("C4").Value2 = "Ready for handoff"
$wsDeDe.Range("D4").Value2 = "e2e"
$wsDeDe.Range("E4").Value2 = "ht"
$wsDeDe.Range("F4").Value2 = "'False"
$wsDeDe.Range("G4").Value2 = "31a51f48-553e-4e6f-bf31-eec0f15bb255.4314296d17028732499ce604519a357cf1aa06d4.de-de.xlf"
$wsDeDe.Range("H4").Value2 = "2016-08-12 04:43:20"
$wsDeDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I4").Value2 = "'"
$wsDeDe.Range("J4").Value2 = "'"
$wsDeDe.Range("K4").Value2 = "0001-01-01 00:00:00"
$wsDeDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L4").Value2 = "'"
$wsDeDe.Range("M4").Value2 = "'True"
$wsDeDe.Range("N4").Value2 = "'"
$wsDeDe.Range("O4").Value2 = "'False"
$wsDeDe.Range("P4").Value2 = "'"

$wsDeDe.Range("A5").Value2 = "ef7eeca8-c9e4-46c0-9841-5970a366c45a.md"
$wsDeDe.Range("B5").Value2 = ".md"
$wsDeDe.Range("C5").Value2 = "Ready for handoff"
$wsDeDe.Range("D5").Value2 = "e2e"
$wsDeDe.Range("E5").Value2 = "ht"
$wsDeDe.Range("F5").Value2 = "'False"
$wsDeDe.Range("G5").Value2 = "ef7eeca8-c9e4-46c0-9841-5970a366c45a.8b6e6573c64d4bec93863aad80d8c624a38d0cef.de-de.xlf"
$wsDeDe.Range("H5").Value2 = "2016-08-12 04:43:20"
$wsDeDe.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I5").Value2 = "'"
$wsDeDe.Range("J5").Value2 = "'"
$wsDeDe.Range("K5").Value2 = "0001-01-01 00:00:00"
$wsDeDe.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L5").Value2 = "'"
$wsDeDe.Range("M5").Value2 = "'True"
$wsDeDe.Range("N5").Value2 = "'"
$wsDeDe.Range("O5").Value2 = "'False"
$wsDeDe.Range("P5").Value2 = "'"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/581f422deee33de35aa3257f53aa5bcd25473dcc/e2e/31a51f48-553e-4e6f-bf31-eec0f15bb255.md", "", "", "31a51f48-553e-4e6f-bf31-eec0f15bb255.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/581f422deee33de35aa3257f53aa5bcd25473dcc/e2e/ef7eeca8-c9e4-46c0-9841-5970a366c45a.md", "", "", "ef7eeca8-c9e4-46c0-9841-5970a366c45a.md")

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P5"))
